# Data cleanup: the "Authors" column (E) for several Systematic Literature
# Review rows had inconsistent internal whitespace in the encoded author
# list. Normalize by widening every run of 2+ spaces (the field separator)
# by two extra spaces, leaving single spaces inside names untouched.

function Add-ExtraSpaces($s) {
    $result = ""
    $i = 0
    $len = $s.Length
    while ($i -lt $len) {
        $ch = $s.Substring($i, 1)
        if ($ch -eq ' ') {
            $runStart = $i
            while ($i -lt $len -and $s.Substring($i, 1) -eq ' ') {
                $i++
            }
            $run = $s.Substring($runStart, $i - $runStart)
            if (($i - $runStart) -ge 2) {
                $run = $run + '  '
            }
            $result = $result + $run
        } else {
            $result = $result + $ch
            $i++
        }
    }
    return $result
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value2 = Add-ExtraSpaces $cell.Value2
}
